$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 187.70589
$ws.Range("I5").Value = 98.90000000000001
$ws.Range("J5").Value = 314.57144
$ws.Range("K5").Value = 98.90000000000001
$ws.Range("L5").Value = 314.57144
$ws.Range("M5").Value = 16.09999999999999
$ws.Range("N5").Value = -544.5714399999999
$ws.Range("H98").Value = 2466.7778
$ws.Range("I98").Value = 2466.7778
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2466.7778
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -968.7777999999998
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 2466.7778
$ws.Range("I122").Value = 2466.7778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7400.3334
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4950.3334
$ws.Range("N122").Value = -4950.3334
$ws.Range("H132").Value = 10757240
$ws.Range("I132").Value = 15875471
$ws.Range("J132").Value = 8954
$ws.Range("K132").Value = 47626413
$ws.Range("L132").Value = 26862
$ws.Range("M132").Value = -47623883
$ws.Range("N132").Value = -31922

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 418.3846
$ws.Range("I97").Value = 420.72726
$ws.Range("K97").Value = 420.72726
$ws.Range("M97").Value = 75.27274
$ws.Range("H102").Value = 41692332
$ws.Range("I102").Value = 41692332
$ws.Range("K102").Value = 41692332
$ws.Range("M102").Value = -41690710
$ws.Range("H110").Value = 946.36584
$ws.Range("I110").Value = 796.2121
$ws.Range("J110").Value = 1565.75
$ws.Range("K110").Value = 796.2121
$ws.Range("L110").Value = 1565.75
$ws.Range("M110").Value = 1248.7879
$ws.Range("N110").Value = -5655.75
$ws.Range("H122").Value = 1266.6666
$ws.Range("I122").Value = 1083.7142
$ws.Range("J122").Value = 1907
$ws.Range("K122").Value = 3251.1426
$ws.Range("L122").Value = 5721
$ws.Range("M122").Value = -801.1425999999997
$ws.Range("N122").Value = -10621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3611.72
$ws.Range("I86").Value = 3789.1052
$ws.Range("J86").Value = 3050
$ws.Range("K86").Value = 3789.1052
$ws.Range("L86").Value = 3050
$ws.Range("M86").Value = -2666.1052
$ws.Range("N86").Value = -5296
$ws.Range("H89").Value = 3611.72
$ws.Range("I89").Value = 3789.1052
$ws.Range("J89").Value = 3050
$ws.Range("K89").Value = 18945.526
$ws.Range("L89").Value = 15250
$ws.Range("M89").Value = -13329.526
$ws.Range("N89").Value = -26482
$ws.Range("H94").Value = 27778840
$ws.Range("I94").Value = 41667292
$ws.Range("K94").Value = 41667292
$ws.Range("M94").Value = -41666841
$ws.Range("H107").Value = 1702.5834
$ws.Range("I107").Value = 1374.1428
$ws.Range("J107").Value = 2162.4
$ws.Range("K107").Value = 1374.1428
$ws.Range("L107").Value = 2162.4
$ws.Range("M107").Value = 545.8571999999999
$ws.Range("N107").Value = -6002.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 409.81818
$ws.Range("J7").Value = 419.8
$ws.Range("L7").Value = 419.8
$ws.Range("N7").Value = -645.8
$ws.Range("H105").Value = 370
$ws.Range("I105").Value = 287.5
$ws.Range("J105").Value = 700
$ws.Range("K105").Value = 287.5
$ws.Range("L105").Value = 700
$ws.Range("M105").Value = 1459.5
$ws.Range("N105").Value = -4194
$ws.Range("H107").Value = 691.8182
$ws.Range("I107").Value = 571.25
$ws.Range("J107").Value = 760.7143
$ws.Range("K107").Value = 571.25
$ws.Range("L107").Value = 760.7143
$ws.Range("M107").Value = 1348.75
$ws.Range("N107").Value = -4600.7143
$ws.Range("H132").Value = 4555.1714
$ws.Range("I132").Value = 4885.346
$ws.Range("K132").Value = 14656.038
$ws.Range("M132").Value = -12126.038

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 16.958334
$ws.Range("I2").Value = 10.944445
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 65.66667
$ws.Range("L2").Value = 210
$ws.Range("M2").Value = 47.33333
$ws.Range("N2").Value = -436
$ws.Range("H4").Value = 365465.75
$ws.Range("I4").Value = 52775.74
$ws.Range("J4").Value = 695527.4399999999
$ws.Range("K4").Value = 158327.22
$ws.Range("L4").Value = 2086582.32
$ws.Range("M4").Value = -158215.22
$ws.Range("N4").Value = -2086806.32
$ws.Range("H16").Value = 3000
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 9000
$ws.Range("N16").Value = -9346
$ws.Range("H40").Value = 201.14285
$ws.Range("I40").Value = 201
$ws.Range("J40").Value = 202
$ws.Range("K40").Value = 804
$ws.Range("L40").Value = 808
$ws.Range("M40").Value = -735
$ws.Range("N40").Value = -946
$ws.Range("H68").Value = 1357.0312
$ws.Range("I68").Value = 691.7857
$ws.Range("J68").Value = 1874.4445
$ws.Range("K68").Value = 2075.3571
$ws.Range("L68").Value = 5623.333500000001
$ws.Range("M68").Value = -1264.3571
$ws.Range("N68").Value = -7245.333500000001
$ws.Range("H71").Value = 1357.0312
$ws.Range("I71").Value = 691.7857
$ws.Range("J71").Value = 1874.4445
$ws.Range("K71").Value = 6226.071300000001
$ws.Range("L71").Value = 16870.0005
$ws.Range("M71").Value = -2170.071300000001
$ws.Range("N71").Value = -24982.0005
$ws.Range("H75").Value = 1300
$ws.Range("J75").Value = 1300
$ws.Range("L75").Value = 3900
$ws.Range("N75").Value = -5896
$ws.Range("H78").Value = 1300
$ws.Range("J78").Value = 1300
$ws.Range("L78").Value = 11700
$ws.Range("N78").Value = -21684
$ws.Range("H87").Value = 3174.6667
$ws.Range("I87").Value = 2333.3333
$ws.Range("J87").Value = 4016
$ws.Range("K87").Value = 6999.999899999999
$ws.Range("L87").Value = 12048
$ws.Range("M87").Value = -5751.999899999999
$ws.Range("N87").Value = -14544
$ws.Range("H90").Value = 3174.6667
$ws.Range("I90").Value = 2333.3333
$ws.Range("J90").Value = 4016
$ws.Range("K90").Value = 20999.9997
$ws.Range("L90").Value = 36144
$ws.Range("M90").Value = -14759.9997
$ws.Range("N90").Value = -48624
$ws.Range("H103").Value = 2892.3845
$ws.Range("J103").Value = 3991
$ws.Range("L103").Value = 11973
$ws.Range("N103").Value = -13731
$ws.Range("H109").Value = 79515.92
$ws.Range("I109").Value = 143958.14
$ws.Range("J109").Value = 4333.3335
$ws.Range("K109").Value = 431874.42
$ws.Range("L109").Value = 13000.0005
$ws.Range("M109").Value = -430834.42
$ws.Range("N109").Value = -15080.0005
$ws.Range("H114").Value = 668.25
$ws.Range("I114").Value = 345.91666
$ws.Range("J114").Value = 1151.75
$ws.Range("K114").Value = 1037.74998
$ws.Range("L114").Value = 3455.25
$ws.Range("M114").Value = 2216.25002
$ws.Range("N114").Value = -9963.25
$ws.Range("H136").Value = 1811.8
$ws.Range("I136").Value = 1131.5
$ws.Range("K136").Value = 3394.5
$ws.Range("M136").Value = 1705.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32148158
$ws.Range("I70").Value = 62504276
$ws.Range("J70").Value = 20005710
$ws.Range("K70").Value = 62504276
$ws.Range("L70").Value = 20005710
$ws.Range("M70").Value = -62504006
$ws.Range("N70").Value = -20006250
$ws.Range("H73").Value = 32148158
$ws.Range("I73").Value = 62504276
$ws.Range("J73").Value = 20005710
$ws.Range("K73").Value = 62504276
$ws.Range("L73").Value = 20005710
$ws.Range("M73").Value = -62503340
$ws.Range("N73").Value = -20007582
$ws.Range("H102").Value = 1212.4
$ws.Range("I102").Value = 1193.6
$ws.Range("J102").Value = 1259.4
$ws.Range("K102").Value = 1193.6
$ws.Range("L102").Value = 1259.4
$ws.Range("M102").Value = 428.4000000000001
$ws.Range("N102").Value = -4503.4
$ws.Range("H132").Value = 2414.7742
$ws.Range("I132").Value = 2008.6842
$ws.Range("J132").Value = 3057.75
$ws.Range("K132").Value = 6026.0526
$ws.Range("L132").Value = 9173.25
$ws.Range("M132").Value = -3496.0526
$ws.Range("N132").Value = -14233.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 27779990
$ws.Range("I122").Value = 31252114
$ws.Range("K122").Value = 93756342
$ws.Range("M122").Value = -93753892
$ws.Range("H132").Value = 55313.21
$ws.Range("I132").Value = 2499.75
$ws.Range("K132").Value = 7499.25
$ws.Range("M132").Value = -4969.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2500.5
$ws.Range("I81").Value = 2500.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5001
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3940
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2500.5
$ws.Range("I84").Value = 2500.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 25005
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -19701
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 667.3077
$ws.Range("I113").Value = 330.22223
$ws.Range("J113").Value = 1425.75
$ws.Range("K113").Value = 990.66669
$ws.Range("L113").Value = 4277.25
$ws.Range("M113").Value = 1179.33331
$ws.Range("N113").Value = -8617.25
$ws.Range("H126").Value = 55556684
$ws.Range("I126").Value = 83333864
$ws.Range("K126").Value = 250001592
$ws.Range("M126").Value = -249999122
